# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 4fcd48a8... source file row (row 2) on the zh-cn and
# de-de target sheets, and mirror the de-de handoff date onto the Overview
# sheet's "Latest HO Xliff Generate Date" column for the e19dce4f... row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 2 (4fcd48a8-4df8-4ff2-b467-2efc5fedda14.md)
$wsZhCn.Range("H2").Value = "2016-10-19 17:03:04"
$wsZhCn.Range("K2").Value = "2016-10-19 17:03:46"

# de-de sheet, row 2 (4fcd48a8-4df8-4ff2-b467-2efc5fedda14.md)
$wsDeDe.Range("H2").Value = "2016-10-19 17:03:16"
$wsDeDe.Range("K2").Value = "2016-10-19 17:04:05"

# Overview sheet, row 3 (e19dce4f-c2c6-4eb7-8d02-38809d9c4184.md) - Latest HO Xliff Generate Date
$wsOverview.Range("G3").Value = "2016-10-19 17:03:16"
